# Applies the commit:
# - Adds a new "Player Info" sheet as the first sheet with player metadata
# - Renames "MATCH_CARD_LINK" column to "MATCH_CODE" in "ODI Batting" and
#   "ODI Bowling" sheets, and replaces the full match-card URL value with
#   just the numeric match code extracted from that URL.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Player Info" worksheet before the first existing sheet ---
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

# Match the page margins used by the other sheets in this workbook
# (0.75in/0.75in/1in/1in/0.5in/0.5in == 54/54/72/72/36/36 points)
$ps = $playerInfo.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row - force ID to be stored as text (matches inlineStr in target), then
# drop the auto-generated "quote prefix" style so the cell keeps the default style.
$playerInfo.Range("A2").Value = "'6622"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Matthew James Potts"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# Copy the header style (bold, centered, bordered) from the existing sheets
# onto the new header row so it matches the look of the other header rows.
$battingHeaderStyle = $wb.Worksheets.Item("ODI Batting").Range("A1")
$battingHeaderStyle.Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# --- 2. Rename MATCH_CARD_LINK -> MATCH_CODE and store just the match code ---
$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")

$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").Value = "'4619"
$batting.Range("D2").Style = "Normal"

$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = "'4619"
$bowling.Range("B2").Style = "Normal"
